$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's data as a new row (row 54), following the same
# pattern as the existing data rows (A = date serial, B/C/D = win counts).
$ws.Range("A54").Value = 46003
$ws.Range("A54").Style = $ws.Range("A53").Style
$ws.Range("A54").NumberFormat = $ws.Range("A53").NumberFormat

$ws.Range("B54").Value = 119
$ws.Range("C54").Value = 131
$ws.Range("D54").Value = 122
